# Auto-generated edit script applying the Asura_Profits.xlsx diff
# Updates H:N profit-calculation columns on rows across ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 4400.2
$ws.Range("I4").Value = 500.25
$ws.Range("J4").Value = 20000
$ws.Range("K4").Value = 500.25
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = -386.25
$ws.Range("N4").Value = -20228
$ws.Range("H32").Value = 26675.375
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 26675.375
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 26675.375
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -27327.375
$ws.Range("H40").Value = 1460.2778
$ws.Range("I40").Value = 1352.3334
$ws.Range("K40").Value = 1352.3334
$ws.Range("M40").Value = -1177.3334
$ws.Range("H80").Value = 5542.2856
$ws.Range("I80").Value = 223.9
$ws.Range("J80").Value = 10377.182
$ws.Range("K80").Value = 671.7
$ws.Range("L80").Value = 31131.546
$ws.Range("M80").Value = 326.3
$ws.Range("N80").Value = -33127.546
$ws.Range("H83").Value = 5542.2856
$ws.Range("I83").Value = 223.9
$ws.Range("J83").Value = 10377.182
$ws.Range("K83").Value = 2015.1
$ws.Range("L83").Value = 93394.638
$ws.Range("M83").Value = 2976.9
$ws.Range("N83").Value = -103378.638
$ws.Range("H132").Value = 2252.2812
$ws.Range("I132").Value = 2098.5417
$ws.Range("K132").Value = 6295.625100000001
$ws.Range("M132").Value = -3765.625100000001
$ws.Range("H137").Value = 2011.7358
$ws.Range("I137").Value = 1316.7142
$ws.Range("J137").Value = 2790.16
$ws.Range("K137").Value = 3950.1426
$ws.Range("L137").Value = 8370.48
$ws.Range("M137").Value = -1400.1426
$ws.Range("N137").Value = -13470.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1322.3125
$ws.Range("I2").Value = 1371.6154
$ws.Range("J2").Value = 1108.6666
$ws.Range("K2").Value = 1371.6154
$ws.Range("L2").Value = 1108.6666
$ws.Range("M2").Value = -1258.6154
$ws.Range("N2").Value = -1334.6666
$ws.Range("H60").Value = 37500
$ws.Range("I60").Value = 37500
$ws.Range("K60").Value = 37500
$ws.Range("M60").Value = -36767
$ws.Range("H116").Value = 1322.3125
$ws.Range("I116").Value = 1371.6154
$ws.Range("J116").Value = 1108.6666
$ws.Range("K116").Value = 1371.6154
$ws.Range("L116").Value = 1108.6666
$ws.Range("M116").Value = 922.3846000000001
$ws.Range("N116").Value = -5696.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1322.3125
$ws.Range("I3").Value = 1371.6154
$ws.Range("J3").Value = 1108.6666
$ws.Range("K3").Value = 1371.6154
$ws.Range("L3").Value = 1108.6666
$ws.Range("M3").Value = -1257.6154
$ws.Range("N3").Value = -1336.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 19000
$ws.Range("J26").Value = 19000
$ws.Range("L26").Value = 19000
$ws.Range("N26").Value = -19574
$ws.Range("H31").Value = 1665.2603
$ws.Range("I31").Value = 1707.9
$ws.Range("J31").Value = 1613.5758
$ws.Range("K31").Value = 1707.9
$ws.Range("L31").Value = 1613.5758
$ws.Range("M31").Value = -1412.9
$ws.Range("N31").Value = -2203.5758
$ws.Range("H34").Value = 1665.2603
$ws.Range("I34").Value = 1707.9
$ws.Range("J34").Value = 1613.5758
$ws.Range("K34").Value = 1707.9
$ws.Range("L34").Value = 1613.5758
$ws.Range("M34").Value = -1505.9
$ws.Range("N34").Value = -2017.5758
$ws.Range("H35").Value = 729.8571
$ws.Range("I35").Value = 729.8571
$ws.Range("K35").Value = 729.8571
$ws.Range("M35").Value = -435.8570999999999
$ws.Range("H36").Value = 3000
$ws.Range("J36").Value = 3000
$ws.Range("L36").Value = 3000
$ws.Range("N36").Value = -3776
$ws.Range("H40").Value = 3000
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1039.6177
$ws.Range("I5").Value = 1563.3889
$ws.Range("J5").Value = 450.375
$ws.Range("K5").Value = 4690.1667
$ws.Range("L5").Value = 1351.125
$ws.Range("M5").Value = -4578.1667
$ws.Range("N5").Value = -1575.125
$ws.Range("H12").Value = 878386.56
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 920214.44
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 2760643.32
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -2760989.32
$ws.Range("H68").Value = 821.41
$ws.Range("I68").Value = 656.69354
$ws.Range("J68").Value = 1090.1578
$ws.Range("K68").Value = 1970.08062
$ws.Range("L68").Value = 3270.4734
$ws.Range("M68").Value = -1159.08062
$ws.Range("N68").Value = -4892.4734
$ws.Range("H71").Value = 821.41
$ws.Range("I71").Value = 656.69354
$ws.Range("J71").Value = 1090.1578
$ws.Range("K71").Value = 5910.24186
$ws.Range("L71").Value = 9811.4202
$ws.Range("M71").Value = -1854.24186
$ws.Range("N71").Value = -17923.4202
$ws.Range("H107").Value = 1220.9841
$ws.Range("I107").Value = 1198.9445
$ws.Range("J107").Value = 1250.3704
$ws.Range("K107").Value = 3596.8335
$ws.Range("L107").Value = 3751.1112
$ws.Range("M107").Value = -1676.8335
$ws.Range("N107").Value = -7591.1112
$ws.Range("H135").Value = 1039.6177
$ws.Range("I135").Value = 1563.3889
$ws.Range("J135").Value = 450.375
$ws.Range("K135").Value = 14070.5001
$ws.Range("L135").Value = 4053.375
$ws.Range("M135").Value = -11535.5001
$ws.Range("N135").Value = -9123.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3360.8462
$ws.Range("I16").Value = 857.1429
$ws.Range("J16").Value = 6281.8335
$ws.Range("K16").Value = 857.1429
$ws.Range("L16").Value = 6281.8335
$ws.Range("M16").Value = -687.1429
$ws.Range("N16").Value = -6621.8335
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1555.5555
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1555.5555
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2145.5555
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1555.5555
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1555.5555
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1769.5555
$ws.Range("H46").Value = 911.1539
$ws.Range("I46").Value = 868.3333
$ws.Range("J46").Value = 1007.5
$ws.Range("K46").Value = 868.3333
$ws.Range("L46").Value = 1007.5
$ws.Range("M46").Value = -680.3333
$ws.Range("N46").Value = -1383.5
$ws.Range("H55").Value = 316.42856
$ws.Range("I55").Value = 90.75
$ws.Range("J55").Value = 617.3333
$ws.Range("K55").Value = 90.75
$ws.Range("L55").Value = 617.3333
$ws.Range("M55").Value = 82.25
$ws.Range("N55").Value = -963.3333
$ws.Range("H61").Value = 13392.706
$ws.Range("I61").Value = 15634
$ws.Range("J61").Value = 2933.3333
$ws.Range("K61").Value = 15634
$ws.Range("L61").Value = 2933.3333
$ws.Range("M61").Value = -15432
$ws.Range("N61").Value = -3337.3333
$ws.Range("H82").Value = 2796.875
$ws.Range("I82").Value = 845.5
$ws.Range("J82").Value = 4748.25
$ws.Range("K82").Value = 845.5
$ws.Range("L82").Value = 4748.25
$ws.Range("M82").Value = -484.5
$ws.Range("N82").Value = -5470.25
$ws.Range("H85").Value = 2796.875
$ws.Range("I85").Value = 845.5
$ws.Range("J85").Value = 4748.25
$ws.Range("K85").Value = 845.5
$ws.Range("L85").Value = 4748.25
$ws.Range("M85").Value = 402.5
$ws.Range("N85").Value = -7244.25
$ws.Range("H113").Value = 13392.706
$ws.Range("I113").Value = 15634
$ws.Range("J113").Value = 2933.3333
$ws.Range("K113").Value = 15634
$ws.Range("L113").Value = 2933.3333
$ws.Range("M113").Value = -13464
$ws.Range("N113").Value = -7273.3333
$ws.Range("H122").Value = 34621456
$ws.Range("I122").Value = 62503744
$ws.Range("K122").Value = 187511232
$ws.Range("M122").Value = -187508782

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2094.4524
$ws.Range("I132").Value = 2167.739
$ws.Range("K132").Value = 6503.217000000001
$ws.Range("M132").Value = -3973.217000000001

